# Dragon's Treasure (Version 2) review edits ("Added a few more slots"):
#
#  1. Remove the "Meta description: Experience Dragon's Treasure slot for
#     free with high RTP, Wild and Scatter symbols, and free spins mode
#     with expansion function." paragraph that used to sit right after the
#     H1 title.
#
#  2. At the very end of the document, drop the old
#     "Prompt: Create a fun and exciting feature image ..." image-prompt
#     paragraph and replace it with two new paragraphs reusing the
#     title/meta-description copy: a bold title line ("Play Dragon's
#     Treasure Slot Free | RTP and Bonus Features") followed by an italic
#     blurb line ("Experience Dragon's Treasure slot for free with high
#     RTP, Wild and Scatter symbols, and free spins mode with expansion
#     function."), stripped of their old "Meta description"/"Prompt:"
#     labels.

$d = $word.ActiveDocument

# --- Step 1: delete the whole "Meta description" paragraph. -----------
$metaRng = $d.Content
$foundMeta = $metaRng.Find.Execute(
    "Meta description", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $foundMeta) {
    throw "Could not locate the 'Meta description' paragraph"
}
$metaRng.Expand(4) | Out-Null   # wdParagraph -> whole enclosing paragraph
$metaRng.Delete()

# --- Step 2: delete the "Prompt: ..." paragraph, remembering where it --
#     lived, then splice in the two replacement paragraphs at that spot.
$promptRng = $d.Content
$foundPrompt = $promptRng.Find.Execute(
    "Prompt:", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $foundPrompt) {
    throw "Could not locate the 'Prompt:' paragraph"
}
$promptRng.Expand(4) | Out-Null   # wdParagraph -> whole enclosing paragraph
$insertAt = $promptRng.Start
$promptRng.Delete()

$insertionPoint = $d.Range($insertAt, $insertAt)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$titlePara = '<w:p xmlns:w="' + $wNs + '">' + `
                 '<w:r/>' + `
                 '<w:r><w:rPr><w:b/></w:rPr>' + `
                    '<w:t>Play Dragon''s Treasure Slot Free | RTP and Bonus Features</w:t>' + `
                 '</w:r>' + `
             '</w:p>'
$blurbPara  = '<w:p xmlns:w="' + $wNs + '">' + `
                 '<w:r/>' + `
                 '<w:r><w:rPr><w:i/></w:rPr>' + `
                    '<w:t>Experience Dragon''s Treasure slot for free with high RTP, Wild and Scatter symbols, and free spins mode with expansion function.</w:t>' + `
                 '</w:r>' + `
             '</w:p>'

[void]$insertionPoint.InsertXML($titlePara + $blurbPara)
